$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data per upstream diff
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.950.65'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.514.90'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.41'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.92'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.87%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.213'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.667'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.71'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000310'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.66'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.083.06'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '618.16'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +7.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '70.044.16'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.76'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.94'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.529.68'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.995'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.74'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '106.19'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +12.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.66'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.06'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.06'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.01'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.87'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.22'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.09'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.68%  '
$ws.Range("B31").Value = 'dogwifhat'
$ws.Range("C31").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.12'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +5.93%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.48'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.116'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.08'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.41%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.11'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.11%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.699.33'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '514.66'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.56%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.394'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.58%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.60'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.86'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0783'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.139'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.42%  '
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.90'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.94%  '
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.30%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.76'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.57%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.70'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.36'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.41%  '
